$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 872949.75
$ws.Range("J17").Value = 975529.1
$ws.Range("L17").Value = 2926587.3
$ws.Range("N17").Value = -2926923.3

$ws.Range("H32").Value = 3974.2727
$ws.Range("I32").Value = 4088.8572
$ws.Range("J32").Value = 3773.75
$ws.Range("K32").Value = 4088.8572
$ws.Range("L32").Value = 3773.75
$ws.Range("M32").Value = -3762.8572
$ws.Range("N32").Value = -4425.75

$ws.Range("H33").Value = 958320.4
$ws.Range("J33").Value = 1553.3334
$ws.Range("L33").Value = 1553.3334
$ws.Range("N33").Value = -2011.3334

$ws.Range("H40").Value = 2105.9092
$ws.Range("I40").Value = 941.25
$ws.Range("K40").Value = 941.25
$ws.Range("M40").Value = -766.25

$ws.Range("H92").Value = 1067.6666
$ws.Range("I92").Value = 1151.5
$ws.Range("K92").Value = 1151.5
$ws.Range("M92").Value = 96.5

$ws.Range("H103").Value = 1666.6666
$ws.Range("J103").Value = 1666.6666
$ws.Range("L103").Value = 4999.9998
$ws.Range("N103").Value = -6171.9998

$ws.Range("H136").Value = 238280
$ws.Range("J136").Value = 238280
$ws.Range("L136").Value = 238280
$ws.Range("N136").Value = -248480

$ws.Range("H137").Value = 64226.125
$ws.Range("I137").Value = 1972.6
$ws.Range("J137").Value = 167982
$ws.Range("K137").Value = 5917.799999999999
$ws.Range("L137").Value = 503946
$ws.Range("M137").Value = -3367.799999999999
$ws.Range("N137").Value = -509046

$ws.Range("H138").Value = 2865.4285
$ws.Range("I138").Value = 1451.3077
$ws.Range("J138").Value = 5163.375
$ws.Range("K138").Value = 4353.9231
$ws.Range("L138").Value = 15490.125
$ws.Range("M138").Value = 786.0769
$ws.Range("N138").Value = -25770.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 157604.33
$ws.Range("I32").Value = 155648.69
$ws.Range("K32").Value = 155648.69
$ws.Range("M32").Value = -155361.69

$ws.Range("H45").Value = 18781.158
$ws.Range("I45").Value = 19071.75
$ws.Range("K45").Value = 19071.75
$ws.Range("M45").Value = -18694.75

$ws.Range("H61").Value = 2838.1428
$ws.Range("I61").Value = 2838.1428
$ws.Range("K61").Value = 2838.1428
$ws.Range("M61").Value = -2626.1428

$ws.Range("H74").Value = 1940.5385
$ws.Range("I74").Value = 1278.75
$ws.Range("J74").Value = 2999.4
$ws.Range("K74").Value = 1278.75
$ws.Range("L74").Value = 2999.4
$ws.Range("M74").Value = -404.75
$ws.Range("N74").Value = -4747.4

$ws.Range("H77").Value = 1940.5385
$ws.Range("I77").Value = 1278.75
$ws.Range("J77").Value = 2999.4
$ws.Range("K77").Value = 6393.75
$ws.Range("L77").Value = 14997
$ws.Range("M77").Value = -2025.75
$ws.Range("N77").Value = -23733

$ws.Range("H97").Value = 2441.85
$ws.Range("I97").Value = 1614.25
$ws.Range("J97").Value = 5752.25
$ws.Range("K97").Value = 1614.25
$ws.Range("L97").Value = 5752.25
$ws.Range("M97").Value = -1118.25
$ws.Range("N97").Value = -6744.25

$ws.Range("H110").Value = 1385.2858
$ws.Range("I110").Value = 1417.7273
$ws.Range("K110").Value = 1417.7273
$ws.Range("M110").Value = 627.2727

$ws.Range("H132").Value = 1766.375
$ws.Range("I132").Value = 1747.4286
$ws.Range("K132").Value = 5242.2858
$ws.Range("M132").Value = -2712.2858

$ws.Range("H136").Value = 2838.1428
$ws.Range("I136").Value = 2838.1428
$ws.Range("K136").Value = 8514.428400000001
$ws.Range("M136").Value = -5964.428400000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7694307.5
$ws.Range("I105").Value = 12502241
$ws.Range("K105").Value = 12502241
$ws.Range("M105").Value = -12500494

$ws.Range("H107").Value = 3744.2144
$ws.Range("I107").Value = 4398.273
$ws.Range("K107").Value = 4398.273
$ws.Range("M107").Value = -2478.273

$ws.Range("H134").Value = 1925.7
$ws.Range("I134").Value = 1711.72
$ws.Range("J134").Value = 2995.6
$ws.Range("K134").Value = 5135.16
$ws.Range("L134").Value = 8986.799999999999
$ws.Range("M134").Value = -2600.16
$ws.Range("N134").Value = -14056.8

$ws.Range("H140").Value = 95780
$ws.Range("J140").Value = 95780
$ws.Range("L140").Value = 95780
$ws.Range("N140").Value = -106140

$ws.Range("H141").Value = 87367.375
$ws.Range("J141").Value = 87367.375
$ws.Range("L141").Value = 87367.375
$ws.Range("N141").Value = -97727.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 17499
$ws.Range("I36").Value = 17499
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 17499
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = -17111
$ws.Range("N36").Value = 0

$ws.Range("H40").Value = 17499
$ws.Range("I40").Value = 17499
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 17499
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = -17339
$ws.Range("N40").Value = 0

$ws.Range("H99").Value = 1925.6666
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 1888.5
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 1888.5
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -4884.5

$ws.Range("H126").Value = 1925.6666
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 1888.5
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 5665.5
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -10605.5

$ws.Range("H132").Value = 2666.5
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws.Range("H134").Value = 3329.25
$ws.Range("I134").Value = 2401.5
$ws.Range("K134").Value = 7204.5
$ws.Range("M134").Value = -4669.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.666668
$ws.Range("I2").Value = 25.9
$ws.Range("J2").Value = 132.5
$ws.Range("K2").Value = 155.4
$ws.Range("L2").Value = 795
$ws.Range("M2").Value = -42.39999999999998
$ws.Range("N2").Value = -1021

$ws.Range("H23").Value = 517.6
$ws.Range("J23").Value = 529.0909
$ws.Range("L23").Value = 1587.2727
$ws.Range("N23").Value = -2057.2727

$ws.Range("H41").Value = 202.5
$ws.Range("I41").Value = 202.5
$ws.Range("K41").Value = 607.5
$ws.Range("M41").Value = -269.5

$ws.Range("H137").Value = 3795
$ws.Range("I137").Value = 1791
$ws.Range("J137").Value = 5131
$ws.Range("K137").Value = 5373
$ws.Range("L137").Value = 15393
$ws.Range("M137").Value = -273
$ws.Range("N137").Value = -25593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10443.777
$ws.Range("I70").Value = 11124.875
$ws.Range("J70").Value = 4995
$ws.Range("K70").Value = 11124.875
$ws.Range("L70").Value = 4995
$ws.Range("M70").Value = -10854.875
$ws.Range("N70").Value = -5535

$ws.Range("H73").Value = 10443.777
$ws.Range("I73").Value = 11124.875
$ws.Range("J73").Value = 4995
$ws.Range("K73").Value = 11124.875
$ws.Range("L73").Value = 4995
$ws.Range("M73").Value = -10188.875
$ws.Range("N73").Value = -6867

$ws.Range("H102").Value = 624
$ws.Range("I102").Value = 582.2222
$ws.Range("K102").Value = 582.2222
$ws.Range("M102").Value = 1039.7778

$ws.Range("H122").Value = 5187.857
$ws.Range("I122").Value = 4802.5
$ws.Range("K122").Value = 14407.5
$ws.Range("M122").Value = -11957.5

$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 27809318
$ws.Range("I61").Value = 37041240
$ws.Range("J61").Value = 113552.664
$ws.Range("K61").Value = 37041240
$ws.Range("L61").Value = 113552.664
$ws.Range("M61").Value = -37041038
$ws.Range("N61").Value = -113956.664

$ws.Range("H100").Value = 32102.637
$ws.Range("I100").Value = 2099.5557
$ws.Range("K100").Value = 2099.5557
$ws.Range("M100").Value = -1558.5557

$ws.Range("H113").Value = 27809318
$ws.Range("I113").Value = 37041240
$ws.Range("J113").Value = 113552.664
$ws.Range("K113").Value = 37041240
$ws.Range("L113").Value = 113552.664
$ws.Range("M113").Value = -37039070
$ws.Range("N113").Value = -117892.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2292
$ws.Range("I81").Value = 2365.125
$ws.Range("K81").Value = 4730.25
$ws.Range("M81").Value = -3669.25

$ws.Range("H84").Value = 2292
$ws.Range("I84").Value = 2365.125
$ws.Range("K84").Value = 23651.25
$ws.Range("M84").Value = -18347.25

$ws.Range("H107").Value = 29412224
$ws.Range("I107").Value = 441.2
$ws.Range("K107").Value = 1323.6
$ws.Range("M107").Value = 596.4000000000001

$ws.Range("H132").Value = 13711.777
$ws.Range("I132").Value = 16501
$ws.Range("K132").Value = 49503
$ws.Range("M132").Value = -46973

$ws.Range("H136").Value = 1506.6666
$ws.Range("I136").Value = 1597.4667
$ws.Range("J136").Value = 1279.6666
$ws.Range("K136").Value = 4792.4001
$ws.Range("L136").Value = 3838.9998
$ws.Range("M136").Value = -2242.4001
$ws.Range("N136").Value = -8938.9998
